$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.445722818374634
$ws.Range("B1").Value = 3.542949914932251
$ws.Range("C1").Value = 2.125890970230103
$ws.Range("D1").Value = 1.124770760536194
$ws.Range("E1").Value = 0.7687835693359375
